$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Header line of the "New in this update" section
Replace-Text "New in this update (Render DB env hardening)" "New in this update (Branding cleanup)"

# First bullet under the section
Replace-Text "- Updated backend DB connection resolver in ``Program.cs`` to support multiple env keys:" "- Removed Emergent branding artifacts from frontend public entry:"

# Sub-bullets: first four are rewritten in place
Replace-Text "  - ``POSTGRES_CONNECTION_STRING``" "  - Removed Emergent badge block (``Made with Emergent``)."
Replace-Text "  - ``ConnectionStrings__Postgres`` / ``ConnectionStrings:Postgres``" "  - Removed Emergent external scripts from ``index.html``."
Replace-Text "  - ``DATABASE_URL``" "  - Updated page title to ``Sitesellr``."
Replace-Text "  - ``RENDER_EXTERNAL_DATABASE_URL``" "  - Updated meta description to ``Sitesellr commerce platform``."

# Remaining old paragraphs that no longer have a counterpart are removed entirely.
# Delete from the bottom up so earlier paragraph indices stay valid.
$targets = @(
    "  - ``RENDER_INTERNAL_DATABASE_URL``",
    "- Added support to normalize ``postgres://...`` URLs into Npgsql connection string format.",
    "- Added production safeguard: if resolved host is localhost, app fails fast with clear error message."
)

foreach ($target in $targets) {
    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd("`r", "`a") -eq $target) {
            $p.Range.Delete()
            break
        }
    }
}

# Git state section updates
Replace-Text "- Last pushed commit: 0c88abb" "- Last pushed commit: 9607153"
Replace-Text "- Current DB env resilience fix is local and not pushed yet." "- Current branding cleanup is local and not pushed yet."
